$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.03134791666666667
$ws.Range("B2").Value = 27.48567361111111
$ws.Range("B3").Value = 1392.685211088334

$ws.Range("B6").Value = 97
$ws.Range("B7").Value = 28

$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 28.99021316195733

$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 48.03984031810768

$ws.Range("A10").Value = "Total SOC consumed(%)"

$ws.Range("B11").Value = "Custom mode`n93.26%`nEco mode`n3.20%`nSports mode`n0.11%"

$ws.Range("A12").Value = "Peak Power(kW)"

$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -1853.481242113863

$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B14").Value = 0.1270907838888889

$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.009124760293835223

$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.33

$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.06

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("B18").Value = 0.27

$ws.Range("A19").Value = "Minimum Temperature(C)"

$ws.Range("A20").Value = "Maximum Temperature(C)"

$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 17

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

$ws.Range("A28").Value = "highest cell temp(C)"

$ws.Range("A29").Value = "lowest cell temp(C)"

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.456740701388889

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.000000149427693807329

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 15.56871943704528

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 5.275712638651453

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 6.687075100385639

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 9.859658887607841

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 10.7979167494931

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 9.016817079473622

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 15.3142766270425

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 14.79346400031806

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 12.55516240607482

$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
